$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Test Status" column (AH) with a header and per-row PASS/FAIL results.
# Row 4 (TEST_AUTO_OAR_08032022_03) failed because no PR was found.
$ws.Range("AH1").Value = "Test Status"

$ws.Range("AH2").Value = "PASS"
$ws.Range("AH2").Interior.ColorIndex = 42

$ws.Range("AH3").Value = "PASS"
$ws.Range("AH3").Interior.ColorIndex = 42

$ws.Range("AH4").Value = "FAIL"
$ws.Range("AH4").Interior.ColorIndex = 10

$ws.Range("AH5").Value = "PASS"
$ws.Range("AH5").Interior.ColorIndex = 42
